$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25 (columns B, D, E, F, G, K, L, M)
# Each inner array: @(row, valB, valD, valE, valF, valG, valK, valL, valM)
$data = @(
    @(2,17.37131022859002,10.60893782362525,14.23731866831932,54.32545083618871,3.745128041384726,11.39029763452691,9.528849336290754,16.34350017344744),
    @(3,17.34399358617459,10.47038548838081,14.17000276143773,53.16764341043818,3.749142674727254,11.35246705026425,9.54859874891388,16.36830337400906),
    @(4,17.33143296835376,10.38359068744915,14.12694154923682,52.44701645988397,3.751731923495542,11.33603206860081,9.562889261545513,16.38675426709716),
    @(5,17.32737760656437,10.34780374863775,14.10895648747819,52.15121351944946,3.752818435168901,11.33104680180436,9.569256127685984,16.3950828316212),
    @(6,17.32676853242085,10.34183661815605,14.10594359854151,52.10197543190171,3.753000748156808,11.3303224561688,9.570346131504191,16.39651466982767),
    @(7,17.33137396682813,10.38310972079998,14.12670077007481,52.44303541901488,3.751746449380605,11.33595790076431,9.562972928513073,16.38686331142128),
    @(8,17.36101997557033,10.56153207646486,14.21446547779806,53.92843731356681,3.746486577035969,11.37584716357758,9.535209247541506,16.35138345939686),
    @(9,17.45235803281999,10.89706548406584,14.37296840180064,56.75133973158523,3.737151952933579,11.50762505230226,9.497972192769723,16.3073834379015),
    @(10,17.53935860291965,11.13395358277016,14.48129029450949,58.7543418614801,3.730882976464295,11.63638406721129,9.481145959012954,16.29065886513044),
    @(11,17.58316510646494,11.23945165641747,14.52883849940475,59.64702767196655,3.728157223183887,11.70167532246904,9.475784556584735,16.2864366906242),
    @(12,17.6003524833137,11.27906109020207,14.54659802341144,59.98216648512742,3.727143038689221,11.72734249918851,9.474084329497421,16.28532421819884),
    @(13,17.59662440058369,11.2705458607071,14.54278409047369,59.91012110361562,3.727360662932223,11.72177308198704,9.474435822881926,16.28554218798725),
    @(14,17.58456717812676,11.24271725386109,14.5303045007923,59.67465923627505,3.728073425599078,11.70376822458386,9.475638062597119,16.286335422936),
    @(15,17.57725947026148,11.22562666584709,14.52262843717275,59.53004719510803,3.728512353633775,11.69286177159819,9.476417453836484,16.2868846233672),
    @(16,17.53658006957713,11.12701220123299,14.47814835250719,58.69560952072575,3.731063636348726,11.63225039156388,9.481542507764301,16.29100288818737),
    @(17,17.51270143741868,11.06592490457324,14.45042048546517,58.17879510646278,3.732660957624689,11.59677229656988,9.485274095696481,16.29439622529985),
    @(18,17.49936579770221,11.03057691322195,14.43430935345041,57.87981025983502,3.733591565078171,11.57700015762377,9.487636219549035,16.2966667263151),
    @(19,17.49491933181183,11.01857270932063,14.42882636365271,57.77829021597762,3.733908695341956,11.57041519046901,9.488473048760095,16.29749023559403),
    @(20,17.51520215817772,11.07244982256538,14.45338899472829,58.2339914206996,3.732489692400343,11.60048355138443,9.484854523338241,16.29400201312514),
    @(21,17.58809250422875,11.25090054462562,14.53397671358508,59.74390070701902,3.727863582428182,11.70903129952133,9.475275977344985,16.28608923576364),
    @(22,17.63921591449793,11.36553779305851,14.58521304628826,60.71369800692695,3.724945007064176,11.78545618175706,9.470939487265371,16.28375234136109),
    @(23,17.61161479432889,11.3045406585677,14.55799747786436,60.19773189551208,3.726493152554053,11.74417350549307,9.473077878791738,16.28474045014004),
    @(24,17.51407035905824,11.06950061645594,14.45204746143082,58.20904294611638,3.732567083130089,11.59880374558,9.485043536790954,16.29417924084799),
    @(25,17.42412633742456,10.80792356508406,14.33152869013622,55.99902398589776,3.739573156806155,11.46630184739796,9.506198814409915,16.31654666638036)
)

$colLetters = @("B","D","E","F","G","K","L","M")

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}

$wb.Save()
